# Auto-generated edit script applying numeric value updates to the
# Kujata_Profits workbook (currentAveragePrice / LevePrice / LeveProfit
# columns), per the scheduled-runner refresh diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 4046.4   # H64
$ws.Cells.Item(64, 9).Value = 4118.8   # I64
$ws.Cells.Item(64, 11).Value = 4118.8   # K64
$ws.Cells.Item(64, 13).Value = -3870.8   # M64
$ws.Cells.Item(67, 8).Value = 4046.4   # H67
$ws.Cells.Item(67, 9).Value = 4118.8   # I67
$ws.Cells.Item(67, 11).Value = 4118.8   # K67
$ws.Cells.Item(67, 13).Value = -3260.8   # M67
$ws.Cells.Item(74, 8).Value = 8000   # H74
$ws.Cells.Item(74, 9).Value = 0   # I74
$ws.Cells.Item(74, 10).Value = 8000   # J74
$ws.Cells.Item(74, 11).Value = 0   # K74
$ws.Cells.Item(74, 12).Value = 8000   # L74
$ws.Cells.Item(74, 13).ClearContents()   # M74
$ws.Cells.Item(74, 14).Value = -9872   # N74
$ws.Cells.Item(76, 8).Value = 6075.25   # H76
$ws.Cells.Item(76, 9).Value = 2650   # I76
$ws.Cells.Item(76, 10).Value = 7217   # J76
$ws.Cells.Item(76, 11).Value = 2650   # K76
$ws.Cells.Item(76, 12).Value = 7217   # L76
$ws.Cells.Item(76, 13).Value = -2335   # M76
$ws.Cells.Item(76, 14).Value = -7847   # N76
$ws.Cells.Item(77, 8).Value = 8000   # H77
$ws.Cells.Item(77, 9).Value = 0   # I77
$ws.Cells.Item(77, 10).Value = 8000   # J77
$ws.Cells.Item(77, 11).Value = 0   # K77
$ws.Cells.Item(77, 12).Value = 40000   # L77
$ws.Cells.Item(77, 13).ClearContents()   # M77
$ws.Cells.Item(77, 14).Value = -49360   # N77
$ws.Cells.Item(79, 8).Value = 6075.25   # H79
$ws.Cells.Item(79, 9).Value = 2650   # I79
$ws.Cells.Item(79, 10).Value = 7217   # J79
$ws.Cells.Item(79, 11).Value = 2650   # K79
$ws.Cells.Item(79, 12).Value = 7217   # L79
$ws.Cells.Item(79, 13).Value = -1558   # M79
$ws.Cells.Item(79, 14).Value = -9401   # N79
$ws.Cells.Item(111, 8).Value = 1692.3334   # H111
$ws.Cells.Item(111, 9).Value = 2447.6667   # I111
$ws.Cells.Item(111, 10).Value = 937   # J111
$ws.Cells.Item(111, 11).Value = 7343.000100000001   # K111
$ws.Cells.Item(111, 12).Value = 2811   # L111
$ws.Cells.Item(111, 13).Value = -4276.000100000001   # M111
$ws.Cells.Item(111, 14).Value = -8945   # N111
$ws.Cells.Item(127, 8).Value = 1783.5   # H127
$ws.Cells.Item(127, 9).Value = 348.33334   # I127
$ws.Cells.Item(127, 10).Value = 2398.5715   # J127
$ws.Cells.Item(127, 11).Value = 1045.00002   # K127
$ws.Cells.Item(127, 12).Value = 7195.7145   # L127
$ws.Cells.Item(127, 13).Value = 3914.99998   # M127
$ws.Cells.Item(127, 14).Value = -17115.7145   # N127
$ws.Cells.Item(138, 8).Value = 1488.356   # H138
$ws.Cells.Item(138, 9).Value = 1274.6061   # I138
$ws.Cells.Item(138, 10).Value = 1759.6538   # J138
$ws.Cells.Item(138, 11).Value = 3823.8183   # K138
$ws.Cells.Item(138, 12).Value = 5278.9614   # L138
$ws.Cells.Item(138, 13).Value = 1316.1817   # M138
$ws.Cells.Item(138, 14).Value = -15558.9614   # N138
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 12000   # H24
$ws.Cells.Item(24, 10).Value = 12000   # J24
$ws.Cells.Item(24, 12).Value = 12000   # L24
$ws.Cells.Item(24, 14).Value = -12748   # N24
$ws.Cells.Item(32, 8).Value = 4058.4119   # H32
$ws.Cells.Item(32, 9).Value = 3694.8914   # I32
$ws.Cells.Item(32, 10).Value = 7402.8   # J32
$ws.Cells.Item(32, 11).Value = 3694.8914   # K32
$ws.Cells.Item(32, 12).Value = 7402.8   # L32
$ws.Cells.Item(32, 13).Value = -3407.8914   # M32
$ws.Cells.Item(32, 14).Value = -7976.8   # N32
$ws.Cells.Item(100, 8).Value = 12000   # H100
$ws.Cells.Item(100, 10).Value = 12000   # J100
$ws.Cells.Item(100, 12).Value = 12000   # L100
$ws.Cells.Item(100, 14).Value = -14164   # N100
$ws.Cells.Item(110, 8).Value = 1302.3334   # H110
$ws.Cells.Item(110, 9).Value = 816.9167   # I110
$ws.Cells.Item(110, 10).Value = 2273.1667   # J110
$ws.Cells.Item(110, 11).Value = 816.9167   # K110
$ws.Cells.Item(110, 12).Value = 2273.1667   # L110
$ws.Cells.Item(110, 13).Value = 1228.0833   # M110
$ws.Cells.Item(110, 14).Value = -6363.1667   # N110
$ws.Cells.Item(132, 8).Value = 1506.2572   # H132
$ws.Cells.Item(132, 9).Value = 1221.7587   # I132
$ws.Cells.Item(132, 10).Value = 2881.3333   # J132
$ws.Cells.Item(132, 11).Value = 3665.2761   # K132
$ws.Cells.Item(132, 12).Value = 8643.999899999999   # L132
$ws.Cells.Item(132, 13).Value = -1135.2761   # M132
$ws.Cells.Item(132, 14).Value = -13703.9999   # N132
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 8337.833000000001   # H134
$ws.Cells.Item(134, 9).Value = 1160.1538   # I134
$ws.Cells.Item(134, 10).Value = 26999.8   # J134
$ws.Cells.Item(134, 11).Value = 3480.4614   # K134
$ws.Cells.Item(134, 12).Value = 80999.39999999999   # L134
$ws.Cells.Item(134, 13).Value = -945.4614000000001   # M134
$ws.Cells.Item(134, 14).Value = -86069.39999999999   # N134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1755.375   # H31
$ws.Cells.Item(31, 9).Value = 1612   # I31
$ws.Cells.Item(31, 10).Value = 1994.3334   # J31
$ws.Cells.Item(31, 11).Value = 1612   # K31
$ws.Cells.Item(31, 12).Value = 1994.3334   # L31
$ws.Cells.Item(31, 13).Value = -1317   # M31
$ws.Cells.Item(31, 14).Value = -2584.3334   # N31
$ws.Cells.Item(34, 8).Value = 1755.375   # H34
$ws.Cells.Item(34, 9).Value = 1612   # I34
$ws.Cells.Item(34, 10).Value = 1994.3334   # J34
$ws.Cells.Item(34, 11).Value = 1612   # K34
$ws.Cells.Item(34, 12).Value = 1994.3334   # L34
$ws.Cells.Item(34, 13).Value = -1410   # M34
$ws.Cells.Item(34, 14).Value = -2398.3334   # N34
$ws.Cells.Item(62, 8).Value = 20002220   # H62
$ws.Cells.Item(62, 9).Value = 2417.5   # I62
$ws.Cells.Item(62, 11).Value = 2417.5   # K62
$ws.Cells.Item(62, 13).Value = -1793.5   # M62
$ws.Cells.Item(65, 8).Value = 20002220   # H65
$ws.Cells.Item(65, 9).Value = 2417.5   # I65
$ws.Cells.Item(65, 11).Value = 12087.5   # K65
$ws.Cells.Item(65, 13).Value = -8967.5   # M65
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 14493966   # H131
$ws.Cells.Item(131, 10).Value = 1241.4777   # J131
$ws.Cells.Item(131, 12).Value = 3724.4331   # L131
$ws.Cells.Item(131, 14).Value = -13804.4331   # N131
$ws.Cells.Item(133, 8).Value = 3171.8   # H133
$ws.Cells.Item(133, 9).Value = 980   # I133
$ws.Cells.Item(133, 11).Value = 2940   # K133
$ws.Cells.Item(133, 13).Value = 2120   # M133
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2331.8572   # H126
$ws.Cells.Item(126, 9).Value = 2354.375   # I126
$ws.Cells.Item(126, 10).Value = 2301.8333   # J126
$ws.Cells.Item(126, 11).Value = 7063.125   # K126
$ws.Cells.Item(126, 12).Value = 6905.499899999999   # L126
$ws.Cells.Item(126, 13).Value = -4593.125   # M126
$ws.Cells.Item(126, 14).Value = -11845.4999   # N126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1793   # H7
$ws.Cells.Item(7, 9).Value = 1699.25   # I7
$ws.Cells.Item(7, 10).Value = 1943   # J7
$ws.Cells.Item(7, 11).Value = 1699.25   # K7
$ws.Cells.Item(7, 12).Value = 1943   # L7
$ws.Cells.Item(7, 13).Value = -1587.25   # M7
$ws.Cells.Item(7, 14).Value = -2167   # N7
$ws.Cells.Item(126, 8).Value = 1793   # H126
$ws.Cells.Item(126, 9).Value = 1699.25   # I126
$ws.Cells.Item(126, 10).Value = 1943   # J126
$ws.Cells.Item(126, 11).Value = 5097.75   # K126
$ws.Cells.Item(126, 12).Value = 5829   # L126
$ws.Cells.Item(126, 13).Value = -2627.75   # M126
$ws.Cells.Item(126, 14).Value = -10769   # N126
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 263.33334   # H107
$ws.Cells.Item(107, 9).Value = 260   # I107
$ws.Cells.Item(107, 10).Value = 300   # J107
$ws.Cells.Item(107, 11).Value = 780   # K107
$ws.Cells.Item(107, 12).Value = 900   # L107
$ws.Cells.Item(107, 13).Value = 1140   # M107
$ws.Cells.Item(107, 14).Value = -4740   # N107
$ws.Cells.Item(132, 8).Value = 3315.125   # H132
$ws.Cells.Item(132, 9).Value = 2696.12   # I132
$ws.Cells.Item(132, 10).Value = 5525.857   # J132
$ws.Cells.Item(132, 11).Value = 8088.36   # K132
$ws.Cells.Item(132, 12).Value = 16577.571   # L132
$ws.Cells.Item(132, 13).Value = -5558.36   # M132
$ws.Cells.Item(132, 14).Value = -21637.571   # N132
$ws.Cells.Item(136, 8).Value = 690.1   # H136
$ws.Cells.Item(136, 9).Value = 607.625   # I136
$ws.Cells.Item(136, 10).Value = 1020   # J136
$ws.Cells.Item(136, 11).Value = 1822.875   # K136
$ws.Cells.Item(136, 12).Value = 3060   # L136
$ws.Cells.Item(136, 13).Value = 727.125   # M136
$ws.Cells.Item(136, 14).Value = -8160   # N136
